# Apply trade #39 close-out update across the workbook.
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1197.66   # Current Capital
$summary.Range("B4").Value = -2.34     # Total P&L $
$summary.Range("B5").Value = -1.2      # Total P&L %
$summary.Range("B6").Value = 39        # Total Trades
$summary.Range("B8").Value = 21        # Losing Trades
$summary.Range("B9").Value = 38.46     # Win Rate %

# --- Strategy Status sheet (MarketMaking row, row 4) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 97.66      # Capital
$status.Range("D4").Value = 39         # Trades
$status.Range("E4").Value = -2.34      # P&L $
$status.Range("F4").Value = -2.34      # P&L %
$status.Range("G4").Value = 38.46      # Win Rate %

# --- All Trades & MarketMaking sheets (trade #39, row 40) ---
$tradeSheets = @("All Trades", "MarketMaking")
foreach ($sheetName in $tradeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("G40").Value = 0.31
    $ws.Range("H40").Value = "CLOSED"
    $ws.Range("I40").Value = -68.3673
    $ws.Range("J40").Value = -0.67
    $ws.Range("K40").Value = 97.66
    $ws.Range("P40").Value = "early_exit"
    $ws.Range("Q40").Value = 3.4
}
